# ---------------------------------------------------------------------------
# Adds a "2022-Q4" worksheet (with fund-holding detail) right after the
# "总计" (summary) sheet and before the existing "2022-Q3" sheet, and
# inserts a corresponding summary row into the "总计" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at
#    the top of the data (row 2), pushing the existing quarters down.
# -----------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Capture existing data rows (2..4) before they get overwritten.
$oldRows = @()
for ($r = 2; $r -le 4; $r++) {
    $oldRows += ,@(
        $summary.Cells.Item($r, 2).Value(),
        $summary.Cells.Item($r, 3).Value(),
        $summary.Cells.Item($r, 4).Value()
    )
}

# Write the captured rows back shifted down by one (row2->3, row3->4, row4->5).
for ($i = 2; $i -ge 0; $i--) {
    $destRow = $i + 3
    $data = $oldRows[$i]

    # Make sure the destination row has the same look (bold/centered index
    # cell in column A) as the other data rows - copy formatting from row 2
    # (which already carries the right style) then overwrite the values.
    $summary.Cells.Item(2, 1).Copy($summary.Cells.Item($destRow, 1))

    $summary.Cells.Item($destRow, 1).Value = $i + 1
    $summary.Cells.Item($destRow, 2).Value = $data[0]
    $summary.Cells.Item($destRow, 3).Value = $data[1]
    $summary.Cells.Item($destRow, 4).Value = $data[2]
}

# Now populate the new row 2 with the 2022-Q4 totals.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 5
$summary.Cells.Item(2, 4).Value = 0.7

# -----------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" worksheet right before "2022-Q3".
# -----------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Force fund-code / numeric-looking text columns to stay text so values
# like "160921" or "8.25" are not silently reinterpreted as numbers.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

# Copy the header / index-column look & feel from the "总计" sheet (style
# used for bold, centered, bordered header & index cells).
$summary.Cells.Item(1, 2).Copy($q4.Range("B1:H1"))
$summary.Cells.Item(2, 1).Copy($q4.Range("A2:A6"))

# Header row.
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Data rows.
$data = @(
    @("160921", "大成多策略混合（LOF）A", "8.25", "87.54", "3.59", "0.2962", 9),
    @("016062", "大成多策略混合（LOF）C", "6.19", "87.54", "3.59", "0.2222", 9),
    @("001898", "易方达大健康主题灵活配置混合", "5.13", "92.14", "3.04", "0.1560", 10),
    @("014121", "大成品质医疗股票A", "0.56", "89.89", "4.54", "0.0254", 9),
    @("014122", "大成品质医疗股票C", "0.08", "89.89", "4.54", "0.0036", 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
}

# Re-activate the summary sheet so the workbook opens on the same tab as
# before the edit.
$summary.Activate()
